$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Disease Ontology source_version: v2024-04-30 -> v2024-05-29
$ws.Range("E3").Value = "v2024-05-29"

# Experimental Factor Ontology source_version: v3.66.0 -> v3.67.0
$ws.Range("E4").Value = "v3.67.0"
